$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell, newValue)
$updates = @{
    "展览" = @(
        @{Cell="F4"; Value=3206},
        @{Cell="F5"; Value=669},
        @{Cell="F6"; Value=567},
        @{Cell="F7"; Value=570},
        @{Cell="F8"; Value=426},
        @{Cell="F9"; Value=148},
        @{Cell="F11"; Value=1410},
        @{Cell="F12"; Value=130},
        @{Cell="F13"; Value=1674},
        @{Cell="F14"; Value=28},
        @{Cell="F16"; Value=634},
        @{Cell="F22"; Value=126},
        @{Cell="F26"; Value=95},
        @{Cell="F27"; Value=4221},
        @{Cell="F28"; Value=19},
        @{Cell="F29"; Value=787},
        @{Cell="F31"; Value=2024},
        @{Cell="F33"; Value=1960}
    )
    "全部类型" = @(
        @{Cell="F4"; Value=3206},
        @{Cell="F5"; Value=669},
        @{Cell="F6"; Value=567},
        @{Cell="F7"; Value=570},
        @{Cell="F9"; Value=426},
        @{Cell="F10"; Value=148},
        @{Cell="F12"; Value=1410},
        @{Cell="F13"; Value=130},
        @{Cell="F14"; Value=1674},
        @{Cell="F15"; Value=28},
        @{Cell="F17"; Value=634},
        @{Cell="F23"; Value=126},
        @{Cell="F27"; Value=95},
        @{Cell="F28"; Value=4221},
        @{Cell="F30"; Value=19},
        @{Cell="F32"; Value=787},
        @{Cell="F34"; Value=2024},
        @{Cell="F36"; Value=1960}
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
